$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the F and G column values (rows 2-9) per the diff
$values = @{
    2 = @{ F = 400;  G = 625 }
    3 = @{ F = 625;  G = 875 }
    4 = @{ F = 875;  G = 1040 }
    5 = @{ F = 1040; G = 1190 }
    6 = @{ F = 1190; G = 1380 }
    7 = @{ F = 1380; G = 1580 }
    8 = @{ F = 1580; G = 1735 }
    9 = @{ F = 1735; G = 2020 }
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row].F
    $ws.Range("G$row").Value = $values[$row].G
}

# Update the active cell selection on the sheet from M8 to G8
$ws.Range("G8").Select() | Out-Null

# Update workbook window view settings (position/size of the workbook window)
$excel.ActiveWindow.WindowState = -4143
$excel.ActiveWindow.Left = 14400
$excel.ActiveWindow.Top = 0
$excel.ActiveWindow.Width = 14400
$excel.ActiveWindow.Height = 7800
